$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.1098014276795042, 1.001597774914103, 3.676440339101017, 1.917404584093044, 1.935410742146401, 46)
    3  = @(-0.04263403788362221, 0.9298606546707162, 2.644551444389527, 1.626207687962865, 1.644018222113035, 45)
    4  = @(0.02349174858992152, 0.8547905997665299, 2.326980562185669, 1.525444381872269, 1.542897175107782, 44)
    5  = @(-0.01048034385804278, 0.9300756200516419, 2.77002958153261, 1.664340584595776, 1.684004220309048, 43)
    6  = @(0.06150569597062733, 0.9485881305707295, 2.681158187928369, 1.6374242541041, 1.656102977579422, 42)
    7  = @(0.01952272370762674, 0.900502877365646, 2.519736357723951, 1.587367744955135, 1.606965805954676, 41)
    8  = @(0.0832317405015639, 0.9720704585887642, 2.678613246716973, 1.636646952374571, 1.655352073199598, 40)
    9  = @(0.04039162028388746, 0.9647072254469649, 2.676716768498583, 1.636067470643733, 1.656949687439857, 39)
    10 = @(0.08220434295181167, 0.9843083280804179, 2.693354164668678, 1.641144163280203, 1.661086189779386, 38)
    11 = @(0.04603842961621714, 0.9311122998615439, 2.722267159893553, 1.649929440883322, 1.672036865348211, 37)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
